# Add VAT and Gross to Stark MOP annual parser (HH sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HH")

# New header cells
$ws.Range("J11").Value = "VAT"
$ws.Range("K11").Value = "Gross"

# Row 12 values
$ws.Range("J12").Value = 13.56
$ws.Range("K12").Value = 136.31

# Row 13 values
$ws.Range("J13").Value = 13.56
$ws.Range("K13").Value = 136.31

# Match the new selection left on the sheet after the edit
$ws.Range("K13").Select()
